$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Variables Listing File name to the new version
$ws.Range("B13").Value = "trend_report_variables_v5.xlsx"

# Update the Variables List Indices lower-right cell reference
$ws.Range("D15").Value = "E238"

# Update the active/selected cell in the sheet view
$ws.Range("D16").Select()
